$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Enigmas. Detectives a domicilio 8. Arte, chanchullos y artimañas" / "Martin, Paul" / "978-84-129217-6-2 "
#   -> "El Cantar de Liébana" / "Peridis" / "978-84-670-7234-1 "
$ws.Range("A2").Value = "`n                        `n                        El Cantar de Liébana`n                        `n                    "
$ws.Range("B2").Value = "`n                        Peridis`n                                            "
$ws.Range("C2").Value = "978-84-670-7234-1 "

# Row 3: "Marienbad eléctrico" / "Vila-Matas, Enrique" / "978-84-322-2578-9 "
#   -> "La hija de la novicia" / "Álvarez, Elena" / "978-84-01-03548-7 "
$ws.Range("A3").Value = "`n                        `n                        La hija de la novicia`n                        `n                    "
$ws.Range("B3").Value = "`n                        Álvarez, Elena`n                                            "
$ws.Range("C3").Value = "978-84-01-03548-7 "

# Row 4: "Paraíso Bacuta" / "Domínguez Rodríguez, Mar" / "978-84-608-9409-4 "
#   -> "El salón dorado" / "Corral, José Luis" / "978-84-1314-409-2 "
$ws.Range("A4").Value = "`n                        `n                        El salón dorado`n                        `n                    "
$ws.Range("B4").Value = "`n                        Corral, José Luis`n                                            "
$ws.Range("C4").Value = "978-84-1314-409-2 "

# Re-fit the rows so the (multi-line) text entry doesn't leave a custom row height behind.
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()
